$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price values that look numeric (e.g. "563.43") as well as
# values using "." as a thousands separator (e.g. "64.878.50"). The source
# data must remain plain text, so force the column's number format to Text
# before writing any values - this prevents Excel from silently converting
# numeric-looking strings into actual numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.878.50"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.134.99"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "563.43"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "147.81"
$ws.Range("E6").Value = "  +6.76%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "3.123.09"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").Value = "6.91"
$ws.Range("E10").Value = "  +13.57%  "
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "35.74"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "3.640.94"
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "64.902.59"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "525.89"
$ws.Range("E18").Value = "  +11.14%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.133.96"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "6.70"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  +5.24%  "
$ws.Range("D23").Value = "7.38"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "12.66"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "78.37"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "8.56"
$ws.Range("E27").Value = "  +14.89%  "
$ws.Range("D28").Value = "2.78"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").Value = "2.11"
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "2.64"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "26.07"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "1.15"
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("D34").Value = "555.71"
$ws.Range("E34").Value = "  +12.10%  "
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").Value = "5.98"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "52.76"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0435"
$ws.Range("E38").Value = "  +8.64%  "
$ws.Range("D39").Value = "0.0810"
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.061.83"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +16.12%  "
$ws.Range("D42").Value = "0.120"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").Value = "8.19"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "0.254"
$ws.Range("E44").Value = "  +7.61%  "
$ws.Range("D45").Value = "2.15"
$ws.Range("E45").Value = "  +8.40%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "24.94"
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "119.31"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "2.09"
$ws.Range("E51").Value = "  +4.95%  "

Write-Host "Updated cryptos list"
